$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.186.26"
$ws.Range("E2").Value = "  +0.74%  "

$ws.Range("D3").Value = "1.569.00"
$ws.Range("E3").Value = "  +0.58%  "

$ws.Range("E4").Value = "  +0.54%  "

$ws.Range("D5").Value = "211.33"
$ws.Range("E5").Value = "  +1.91%  "

$ws.Range("E6").Value = "  +0.70%  "

$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("D8").Value = "22.05"
$ws.Range("E8").Value = "  -0.37%  "

$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("D10").Value = "0.0599"
$ws.Range("E10").Value = "  +0.52%  "

$ws.Range("D11").Value = "0.0868"
$ws.Range("E11").Value = "  +0.99%  "

$ws.Range("D12").Value = "1.794.11"
$ws.Range("E12").Value = "  +0.68%  "

$ws.Range("D13").Value = "1.570.06"
$ws.Range("E13").Value = "  +1.54%  "

$ws.Range("D14").Value = "3.78"
$ws.Range("E14").Value = "  +0.70%  "

$ws.Range("D15").Value = "0.519"
$ws.Range("E15").Value = "  -0.17%  "

$ws.Range("D16").Value = "27.185.59"
$ws.Range("E16").Value = "  +0.79%  "

$ws.Range("D17").Value = "62.28"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "7.44"
$ws.Range("E18").Value = "  +0.84%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "216.25"
$ws.Range("E19").Value = "  -0.44%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0701"
$ws.Range("E20").Value = "  -0.64%  "

$ws.Range("E21").Value = "  +0.41%  "

$ws.Range("D22").Value = "4.14"
$ws.Range("E22").Value = "  +1.00%  "

$ws.Range("D23").Value = "9.23"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").Value = "1.94"
$ws.Range("E24").Value = "  +0.87%  "

$ws.Range("D25").Value = "153.75"
$ws.Range("E25").Value = "  +0.28%  "

$ws.Range("D26").Value = "6.67"
$ws.Range("E26").Value = "  +0.73%  "

$ws.Range("D27").Value = "15.08"
$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("E28").Value = "  +2.44%  "

$ws.Range("E29").Value = "  +0.31%  "

$ws.Range("E30").Value = "  +2.49%  "

$ws.Range("D31").Value = "0.0473"
$ws.Range("E31").Value = "  +1.03%  "

$ws.Range("E32").Value = "  +0.10%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "3.18"
$ws.Range("E33").Value = "  +1.97%  "

$ws.Range("B34").Value = "Maker"
$ws.Range("C34").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D34").Value = "1.451.10"
$ws.Range("E34").Value = "  +1.92%  "

$ws.Range("E35").Value = "  +4.60%  "

$ws.Range("E36").Value = "  +0.16%  "

$ws.Range("E37").Value = "  +1.19%  "

$ws.Range("E38").Value = "  +0.87%  "

$ws.Range("D39").Value = "0.536"
$ws.Range("E39").Value = "  +0.59%  "

$ws.Range("D40").Value = "5.81"
$ws.Range("E40").Value = "  +2.06%  "

$ws.Range("D41").Value = "0.808"
$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("E42").Value = "  +0.37%  "

$ws.Range("D43").Value = "2.34"
$ws.Range("E43").Value = "  +0.48%  "

$ws.Range("D44").Value = "0.997"
$ws.Range("E44").Value = "  -0.16%  "

$ws.Range("D45").Value = "64.50"
$ws.Range("E45").Value = "  -0.56%  "

$ws.Range("E46").Value = "  -0.30%  "

$ws.Range("D47").Value = "1.705.49"
$ws.Range("E47").Value = "  +0.59%  "

$ws.Range("D48").Value = "85.86"
$ws.Range("E48").Value = "  -1.74%  "

$ws.Range("E49").Value = "  +3.68%  "

$ws.Range("D50").Value = "0.0523"
$ws.Range("E50").Value = "  +0.21%  "

$ws.Range("D51").Value = "0.0959"
$ws.Range("E51").Value = "  +0.48%  "
